$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L6").Value = 1392.03
$wsGrupo.Range("L23").Value = "1 de 21"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F6").Value = 1392.03
$wsMensual.Range("F23").Value = 32822.98

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D11").Value = 1392.03
$wsCumplimiento.Range("E11").Value = 4452.41916370549
$wsCumplimiento.Range("F11").Value = 0.2381798456977983

$wsCumplimiento.Range("D15").Value = 32822.98
$wsCumplimiento.Range("E15").Value = 22601.76316613378
$wsCumplimiento.Range("F15").Value = 0.592208066740413
